$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '42.691.15'
$ws.Range('E2').Value = '  -0.29%  '
Set-TextValue 'D3' '2.546.42'
$ws.Range('E3').Value = '  +0.30%  '
Set-TextValue 'D4' '0.998'
$ws.Range('E4').Value = '  -0.01%  '
Set-TextValue 'D5' '318.65'
$ws.Range('E5').Value = '  +4.54%  '
Set-TextValue 'D6' '95.28'
$ws.Range('E6').Value = '  -2.29%  '
Set-TextValue 'D7' '0.579'
$ws.Range('E7').Value = '  +0.39%  '
$ws.Range('E8').Value = '  -0.05%  '
Set-TextValue 'D9' '0.533'
$ws.Range('E9').Value = '  -1.87%  '
Set-TextValue 'D10' '36.47'
$ws.Range('E10').Value = '  +0.02%  '
Set-TextValue 'D11' '0.0815'
$ws.Range('E11').Value = '  -1.08%  '
Set-TextValue 'D12' '7.71'
$ws.Range('E12').Value = '  +1.92%  '
$ws.Range('E13').Value = '  -0.15%  '
Set-TextValue 'D14' '2.937.75'
$ws.Range('E14').Value = '  +0.32%  '
Set-TextValue 'D15' '15.94'
$ws.Range('E15').Value = '  +6.28%  '
Set-TextValue 'D16' '2.553.11'
$ws.Range('E16').Value = '  +0.30%  '
Set-TextValue 'D17' '0.870'
$ws.Range('E17').Value = '  +0.74%  '
Set-TextValue 'D18' '42.738.26'
$ws.Range('E18').Value = '  -0.19%  '
Set-TextValue 'D19' '13.11'
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('E20').Value = '  +1.28%  '
Set-TextValue 'D21' '0.0₃0970'
$ws.Range('E21').Value = '  -1.69%  '
Set-TextValue 'D22' '71.09'
$ws.Range('E22').Value = '  -0.72%  '
Set-TextValue 'D23' '253.23'
$ws.Range('E24').Value = '  +1.54%  '
$ws.Range('E25').Value = '  -2.54%  '
Set-TextValue 'D26' '27.30'
$ws.Range('E26').Value = '  -2.05%  '
Set-TextValue 'D27' '1.00'
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('E28').Value = '  +3.30%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D29' '39.55'
$ws.Range('E29').Value = '  +4.35%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D30' '10.25'
$ws.Range('E30').Value = '  +0.89%  '
Set-TextValue 'D31' '5.97'
$ws.Range('E31').Value = '  -3.15%  '
Set-TextValue 'D32' '156.03'
$ws.Range('E32').Value = '  -0.74%  '
Set-TextValue 'D33' '2.16'
$ws.Range('E33').Value = '  +1.33%  '
Set-TextValue 'D34' '3.36'
$ws.Range('E34').Value = '  +1.80%  '
Set-TextValue 'D35' '19.31'
$ws.Range('E35').Value = '  -1.04%  '
Set-TextValue 'D36' '0.0791'
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('E37').Value = '  +0.04%  '
Set-TextValue 'D38' '0.112'
$ws.Range('E38').Value = '  -2.52%  '
Set-TextValue 'D39' '2.45'
$ws.Range('E39').Value = '  +13.59%  '
$ws.Range('E40').Value = '  -0.22%  '
Set-TextValue 'D41' '24.03'
$ws.Range('E41').Value = '  -3.57%  '
Set-TextValue 'D42' '3.85'
$ws.Range('E42').Value = '  -0.09%  '
Set-TextValue 'D43' '3.37'
$ws.Range('E43').Value = '  -1.13%  '
$ws.Range('E44').Value = '  +0.33%  '
Set-TextValue 'D45' '0.0303'
$ws.Range('E45').Value = '  -0.43%  '
Set-TextValue 'D46' '2.031.28'
$ws.Range('E46').Value = '  -2.80%  '
Set-TextValue 'D47' '84.55'
$ws.Range('E47').Value = '  -2.11%  '
Set-TextValue 'D48' '8.94'
$ws.Range('E48').Value = '  +0.17%  '
Set-TextValue 'D49' '2.790.18'
$ws.Range('E49').Value = '  +0.15%  '
Set-TextValue 'D50' '73.88'
$ws.Range('E50').Value = '  +0.44%  '
Set-TextValue 'D51' '0.191'
$ws.Range('E51').Value = '  -0.32%  '
